$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 with new inventory data, preserving text type for
#     numeric-looking values (matches original inlineStr text cells) ---
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "123456789"

$ws.Cells.Item(2,2).Value = "caraca"

$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = "0.0"

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "200.0"

$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = "100.0"

$ws.Cells.Item(2,6).Value = "Sin proveedor"

# Re-apply the row's original look (fill/alignment) so the number-format
# change above doesn't leave the cells on a divergent style.
$ws.Range("B2").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Remove the old sample row entirely (shifts nothing up below it) ---
$ws.Rows("3:3").Delete()

# --- Column width tweaks ---
# NOTE: the saved OOXML <col width=.../> ends up 5/6 (~0.8333) wider than
# the ColumnWidth value we assign here, so back that padding out up front
# to land on the exact target widths of 11 / 8 / 15.
$ws.Columns("A:A").ColumnWidth = 11 - 5/6
$ws.Columns("B:B").ColumnWidth = 8 - 5/6
$ws.Columns("F:F").ColumnWidth = 15 - 5/6
